# Calendar plan v1.3 added + Assignment updated

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Color constants (OLE_COLOR = 0x00BBGGRR) ---
$green  = 5296274   # FF92D050
$yellow = 65535      # FFFFFF00
$orange = 49407      # FFFFC000
$red    = 255        # FFFF0000

# ------------------------------------------------------------------
# 1) Completion-date (F column) updates for already-finished tasks
# ------------------------------------------------------------------
$ws.Range("F2").Value  = 42295.041666666664
$ws.Range("F3").Value  = 42291.833333333336
$ws.Range("F4").Value  = 42294.854166666664
$ws.Range("F6").Value  = 42294.75
$ws.Range("F6").NumberFormat = "m/d/yy h:mm"
$ws.Range("F8").Value  = 42294.954861111109
$ws.Range("F9").Value  = 42294.954861111109
$ws.Range("F10").Value = 42292.708333333336
$ws.Range("F11").Value = 42295.8125

# ------------------------------------------------------------------
# 2) Renumber the following tasks (shift by one) for the extended row
# ------------------------------------------------------------------
$ws.Range("A7").Value  = 6
$ws.Range("A8").Value  = 7
$ws.Range("A9").Value  = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10

# ------------------------------------------------------------------
# 3) Second block (prototyping pages) gets review deadlines in E
# ------------------------------------------------------------------
$ws.Range("E13").Value = 42295.999305555553
$ws.Range("E14").Value = 42295.999305555553
$ws.Range("E15").Value = 42295.999305497688
$ws.Range("E16").Value = 42295.999305497688
$ws.Range("E17").Value = 42295.999305497688
$ws.Range("E18").Value = 42295.999305497688
$ws.Range("E19").Value = 42295.999305497688
$ws.Range("E20").Value = 42295.999305497688

# ------------------------------------------------------------------
# 4) New task row 21
# ------------------------------------------------------------------
$ws.Range("A21").Value = 11
$ws.Range("B21").Value = "Рецензирование прототипов страниц"
$ws.Range("C21").Value = "Бидзиля, Сорокин"
$ws.Range("E21").Value = 42295.993055555555
$ws.Range("E21").NumberFormat = "m/d/yy h:mm"

# ------------------------------------------------------------------
# 5) Note in column G for the extended deadline (row 7) — set here so
#    shared-string indices line up with the source workbook.
# ------------------------------------------------------------------
$ws.Range("G7").Value = "Продлено до 18.10.15 23:59"

# ------------------------------------------------------------------
# 6) New task row 22
# ------------------------------------------------------------------
$ws.Range("A22").Value = 12
$ws.Range("B22").Value = "Рецензирование материалов второй подгруппы"
$ws.Range("C22").Value = "Заварзин, Руданов"

# ------------------------------------------------------------------
# 7) Status colouring of column A (# cell) for each task
#    Order matters: it determines fill/cellXf creation order so it
#    matches green, yellow, orange, red.
# ------------------------------------------------------------------
# Green = completed on time
$ws.Range("A2").Interior.Color  = $green
$ws.Range("A3").Interior.Color  = $green
$ws.Range("A4").Interior.Color  = $green
$ws.Range("A6").Interior.Color  = $green
$ws.Range("A8").Interior.Color  = $green
$ws.Range("A9").Interior.Color  = $green
$ws.Range("A10").Interior.Color = $green
$ws.Range("A11").Interior.Color = $green

# Yellow = pending review (second block + new row 21)
$ws.Range("A13").Interior.Color = $yellow
$ws.Range("A14").Interior.Color = $yellow
$ws.Range("A15").Interior.Color = $yellow
$ws.Range("A16").Interior.Color = $yellow
$ws.Range("A17").Interior.Color = $yellow
$ws.Range("A18").Interior.Color = $yellow
$ws.Range("A19").Interior.Color = $yellow
$ws.Range("A20").Interior.Color = $yellow
$ws.Range("A21").Interior.Color = $yellow

# Orange = deadline extended
$ws.Range("A7").Interior.Color = $orange

# Red = not completed
$ws.Range("A5").Interior.Color = $red

# ------------------------------------------------------------------
# 8) Update the active selection shown in the sheet view
# ------------------------------------------------------------------
$ws.Range("B16").Select()
